$d = $word.ActiveDocument

# 1. Insert a new "Source Code" paragraph (library(readr)) right after the
#    date paragraph, before the first bookmark.
$datePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Date") {
        $datePara = $p
        break
    }
}

$datePara.Range.InsertParagraphAfter()

$newPara = $datePara.Next()
$newPara.Style = "Source Code"

$r1 = $newPara.Range
$r1.Collapse(1)
$r1.InsertAfter("library")
$r1.Style = "FunctionTok"

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("(readr)")
$r2.Style = "NormalTok"

# 2. Fix the typo "successful complete" -> "successfully complete" in the
#    final paragraph of section 5.
$d.Content.Find.Execute("successful complete", $true, $false, $false, $false, $false,
                         $true, 1, $false, "successfully complete", 2)
